# Adding games to leagues
# This script adds two new worksheets (LEAGUE_RULES and GAMES) to the
# workbook, populates them with data, applies number formats/validations,
# and makes the GAMES sheet the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the two new worksheets after the existing ones, in order.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$leagueRules = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$leagueRules.Name = "LEAGUE_RULES"

$games = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $leagueRules)
$games.Name = "GAMES"

# ---------------------------------------------------------------------
# 2. LEAGUE_RULES sheet content
# ---------------------------------------------------------------------
$leagueRules.Columns.Item(2).ColumnWidth = 14.11
$leagueRules.Columns.Item(3).ColumnWidth = 13.55
$leagueRules.Columns.Item(4).ColumnWidth = 14.22
$leagueRules.Columns.Item(5).ColumnWidth = 16.66

# Build the dropdown source/lookup lists first (rows 16-18) ...
$leagueRules.Range("A16").Value = "RUNS"
$leagueRules.Range("A17").Value = "WICKETS"
$leagueRules.Range("A18").Value = "CATCHES"

# ... then the header row label for column A ...
$leagueRules.Range("A1").Value = "METRIC"

# ... the OPERATOR lookup values ...
$leagueRules.Range("B16").Value = "LESS_THAN"
$leagueRules.Range("B17").Value = "GREATER_THAN"
$leagueRules.Range("B18").Value = "EQUALS"

# ... the RATIO lookup values ...
$leagueRules.Range("E16").Value = "PER_SCORE"
$leagueRules.Range("E17").Value = "TOTAL"

# ... then the rest of the header row ...
$leagueRules.Range("B1").Value = "OPERATOR"
$leagueRules.Range("C1").Value = "SCORE"
$leagueRules.Range("D1").Value = "POINTS"
$leagueRules.Range("E1").Value = "RATIO"

# Rule rows (re-use values already present in the shared string table)
$leagueRules.Range("A2").Value = "RUNS"
$leagueRules.Range("B2").Value = "EQUALS"
$leagueRules.Range("C2").Value = 0
$leagueRules.Range("D2").Value = -5
$leagueRules.Range("E2").Value = "TOTAL"

$leagueRules.Range("A3").Value = "RUNS"
$leagueRules.Range("B3").Value = "GREATER_THAN"
$leagueRules.Range("C3").Value = 50
$leagueRules.Range("D3").Value = 5
$leagueRules.Range("E3").Value = "PER_SCORE"

# C4 is an empty, quote-prefixed cell (user typed ' then cleared it)
$leagueRules.Range("C4").Value = "'"
$leagueRules.Range("C4").Value = ""

# Data validation (drop-down lists) driven by the lookup lists above
$leagueRules.Range("A1:A14").Validation.Add(3, 1, 1, "=`$A`$16:`$A`$18")
$leagueRules.Range("B2:B15").Validation.Add(3, 1, 1, "=`$B`$16:`$B`$18")
$leagueRules.Range("E2:E15").Validation.Add(3, 1, 1, "=`$E`$16:`$E`$17")

# ---------------------------------------------------------------------
# 3. GAMES sheet content
# ---------------------------------------------------------------------
$games.Columns.Item(3).ColumnWidth = 19.89
$games.Columns.Item(4).ColumnWidth = 18.11

# Header row (TEAM1 / TEAM2 first) ...
$games.Range("A1").Value = "TEAM1"
$games.Range("B1").Value = "TEAM2"

# ... first game's teams ...
$games.Range("A2").Value = "RR"
$games.Range("B2").Value = "MI"

# ... rest of the header row ...
$games.Range("C1").Value = "DATE( MM/DD/YYYY)"
$games.Range("D1").Value = "TIME(hh:mm:ss)"
$games.Range("E1").Value = "VENUE"

# ... first game's date/time/venue ...
$games.Range("C2").Value = 43066
$games.Range("C2").NumberFormat = "mm-dd-yy"
$games.Range("D2").Value = 0.66666666666666663
$games.Range("D2").NumberFormat = "h:mm:ss"
$games.Range("E2").Value = "MUMBAI"

# Second game (re-uses values already present in the shared string table)
$games.Range("A3").Value = "CSK"
$games.Range("B3").Value = "RCB"
$games.Range("C3").Value = 43066
$games.Range("D3").Value = 0.83333333333333337
$games.Range("E3").Value = "MUMBAI"

# Re-use the exact same date/time cell formats for the second game row
# (copy formats only, so the underlying style record is shared instead of
# a new, functionally-identical one being allocated).
$games.Range("C2").Copy()
$games.Range("C3").PasteSpecial(-4122)
$games.Range("D2").Copy()
$games.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Make GAMES the active/selected sheet & set selections to match
# ---------------------------------------------------------------------
$null = $leagueRules.Range("A4").Select()
$null = $games.Range("D7").Select()
$games.Activate()
